# Cadastro da relacao Usuario X Equipe
# Adds the "tb_c_equipe_usua" table and "sq_equipe_usua" sequence grant
# blocks to the "grants por usuario" sheet, replicating the existing
# per-user grant-row pattern (rows 149-184).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Same 18-user order already used for every other grant block on this sheet.
$users = @(
  'aline',
  '"adriana.paes"',
  '"kleverson.antonio"',
  '"luana.mourao"',
  '"mariana.brider"',
  '"adriana.paes"',
  '"larissa.rocha"',
  '"livia.hallack"',
  '"vanessa.cirilo"',
  '"milena.ribeiral"',
  '"tatiane.neto"',
  '"tatyene.nehrer"',
  '"victor.quinet"',
  '"lidice.lenz"',
  '"alberlania.muller"',
  '"marcia.costa"',
  'vanessa',
  'aline'
)

$tableGrant = "GRANT DELETE, INSERT, SELECT, UPDATE ON TABLE tratamento.tb_c_equipe_usua TO"
$seqGrant = "GRANT ALL ON SEQUENCE tratamento.sq_equipe_usua TO"

$firstRow = 149
$tableLastRow = $firstRow + $users.Length - 1          # 166
$seqLastRow = $tableLastRow + $users.Length             # 184

# --- Table grant rows (149-166) ---
$row = $firstRow
foreach ($u in $users) {
  $ws.Cells.Item($row, 1).Value = $tableGrant
  $ws.Cells.Item($row, 2).Value = $u
  $ws.Cells.Item($row, 3).Value = ";"
  $row++
}

# --- Sequence grant rows (167-184) ---
foreach ($u in $users) {
  $ws.Cells.Item($row, 1).Value = $seqGrant
  $ws.Cells.Item($row, 2).Value = $u
  $ws.Cells.Item($row, 3).Value = ";"
  $row++
}

# Match column B formatting (wrapped, vertical-centered Arial 9) used by
# every other data row, by copying the format from the row right above.
$ws.Range("B130").Copy() | Out-Null
$ws.Range("B$firstRow`:B$seqLastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Concatenation formula for column D, same as the rest of the sheet.
$ws.Range("D$firstRow`:D$seqLastRow").Formula = "=A$firstRow&"" ""&B$firstRow&"" ""&C$firstRow"

# Reset the stale selection/scroll position left over from before the new
# rows existed.
$ws.Range("A1").Select()
